# Automatische test-sync: 2025-06-26 20:55:50
# Adds a new logged e-mail (row 13) to the "Logs" sheet and updates the
# related conditional-formatting ranges + the "Dashboard" summary count.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append the new row of data --------------------------------------
$logs.Range("A13").Value = "Wil je 4 rollen plasticfolie bestellen?"
$logs.Range("B13").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C13").Value = "Hoi Johan,`nWil je 4 rollen plasticfolie bestellen?`nMarc`nSent using {0}"
$logs.Range("D13").Value = "Bestelling / Levering"
$logs.Range("E13").Value = "Beste Marc,`nBedankt voor je e-mail. Helaas kunnen we geen bestellingen via e-mail verwerken. Je kunt eenvoudig via onze website 4 rollen plasticfolie bestellen. Mocht je hulp nodig hebben of vragen hebben over het bestelproces, dan helpen we je graag verder.`nMet vriendelijke groet,`nJohan"
$logs.Range("F13").Value = "2025-06-26 20:54:55"
$logs.Range("G13").Value = "Ja"
$logs.Range("H13").Value = "Nee"
$logs.Range("I13").Value = "Ja"

# Setting multi-line values can mark the row with an explicit custom
# height; AutoFit puts it back into the default "auto" state so the row
# serializes the same way as the pre-existing rows.
$logs.Rows.Item(13).AutoFit()

# --- Extend the conditional formatting ranges to cover row 13 --------
$colRanges = @("D2:D12", "G2:G12", "H2:H12", "I2:I12")
$newLast = @{ "D2:D12" = "D2:D13"; "G2:G12" = "G2:G13"; "H2:H12" = "H2:H13"; "I2:I12" = "I2:I13" }

foreach ($old in $colRanges) {
    $fcs = $logs.Range($old).FormatConditions
    $target = $logs.Range($newLast[$old])
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($target)
    }
}

# --- Update the Dashboard summary count for "Bestelling / Levering" --
$dash.Range("B2").Value = 8
